# Update column F (dSF) values for rows 3-21 (excluding 2, 10, 13, 16 which are unchanged)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 3
    4  = -1
    5  = 1
    6  = 4
    7  = -3
    8  = -4
    9  = -5
    11 = -2
    12 = 2
    14 = 2
    15 = -1
    17 = 2
    18 = 7
    19 = -5
    20 = -4
    21 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
